# FINFLUX-3612 Cartias specific scenarios
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet: update late-fee (penalty) figures for the 3rd installment
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Activate()
$wsSummary.Range("F3").Value = 97.04
$wsSummary.Range("F5").Value = 0.14
$wsSummary.Range("B9").Select()

# ---------------------------------------------------------------------------
# Repayment schedule sheet: just move the cursor (cosmetic selection change)
# ---------------------------------------------------------------------------
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Activate()
$wsRepay.Range("J9").Select()

# ---------------------------------------------------------------------------
# Transactions sheet: renumber the transaction IDs in column A
# ---------------------------------------------------------------------------
$wsTxns = $wb.Worksheets.Item("Transactions")
$wsTxns.Activate()
$wsTxns.Range("A2").Value = 427
$wsTxns.Range("A3").Value = 426
$wsTxns.Range("A4").Value = 441
$wsTxns.Range("A5").Value = 440
$wsTxns.Range("A6").Value = 434
$wsTxns.Range("A7").Value = 424
$wsTxns.Range("A8").Value = 423
$wsTxns.Range("A9").Value = 422
$wsTxns.Range("A10").Value = 421
$wsTxns.Range("I4").Select()

# ---------------------------------------------------------------------------
# ChargesTab sheet: widen column F to match column E, shrink rows 3 & 4,
# and point the waive-charge hyperlinks at the new loan account charge ids.
# ---------------------------------------------------------------------------
$wsCharges = $wb.Worksheets.Item("ChargesTab")
$wsCharges.Activate()

$wsCharges.Columns.Item(6).ColumnWidth = $wsCharges.Columns.Item(5).ColumnWidth

$wsCharges.Rows.Item(3).RowHeight = 45
$wsCharges.Rows.Item(4).RowHeight = 45

$hyperlinkAddress = "file:///D:\bharath\Code\conflux-web\app\index.html%3fbaseApiUrl=https:\localhost:8443\fineract-provider\api\v1&tenantIdentifier=default"
$displayBase = "D:\bharath\Code\conflux-web\app\index.html?baseApiUrl=https:\localhost:8443\fineract-provider\api\v1&tenantIdentifier=default"

$wsCharges.Hyperlinks.Delete()
$wsCharges.Hyperlinks.Add($wsCharges.Range("K3"), $hyperlinkAddress, "/loanaccountcharge/20/waivecharge/88", "", "$displayBase - /loanaccountcharge/20/waivecharge/88")
$wsCharges.Hyperlinks.Add($wsCharges.Range("K4"), $hyperlinkAddress, "/loanaccountcharge/20/waivecharge/91", "", "$displayBase - /loanaccountcharge/20/waivecharge/91")

$wsCharges.Range("G6").Select()
